$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo'd role names in the header row (row 1)
$ws.Range("BG1").Value = "WM Microsoft Project"
$ws.Range("BJ1").Value = "AP Invoice Processor  (GFT Job Role)"

# Reflect the saved view's active cell in the frozen (bottom-right) pane
$ws.Activate()
$ws.Range("BJ1").Select()
